# Re-apply the new matrices ranking (mat_rank) results after excluding
# the "no rank decision" entries from the binary classification.
# This reorders the twelve female workers by their recomputed matrices
# score (column F) and refreshes each worker's new score, while the
# worker's own attributes (id, name, gender, race) travel with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# prolificid, B (numeric id), name, gender, F score, race, H (rank)
$rows = @(
    @("5c5882fc5bfe7600011197cb", 2,  "Colleen",   "female", 13.42119510329043, "White", 1),
    @("60bd88b8fc436774352f53b9", 3,  "Annes",     "female", 13.17773416771519, "Asian", 2),
    @("608b14a312c099ac00b721b6", 22, "Khushi",    "female", 8.277947983434146, "Asian", 3),
    @("5c0e89c6c323400001e6c4a5", 21, "Bri",       "female", 8.218874334828817, "Black or African American", 4),
    @("60b45e9961dd412bfb6780f8", 19, "Jewel",     "female", 8.21192345112825,  "Black or African American", 5),
    @("60cb36ee9f58331a33cf5506", 33, "Shaniek",   "female", 5.441970684512863, "Black or African American", 6),
    @("6036f9b3b1842f8b659b18c7", 32, "Kellie",    "female", 5.381459162249058, "White", 7),
    @("60d5775a99b502eec8cf56b4", 30, "Shadaisia", "female", 5.321845954194636, "Black or African American", 8),
    @("5e96194b0a9fe909389e9f7b", 34, "Tina",      "female", 4.498467056693604, "White", 9),
    @("6077db0613ce87b4a62a78f9", 35, "Lori",      "female", 4.222996349665409, "White", 10),
    @("60bfcf5805c5ae12a546f9f3", 41, "Giana",     "female", 2.390791975163696, "White", 11),
    @("60c0e5899d387663c07eb3a4", 44, "Nansi",     "female", 1.089220531548616, "Asian", 12)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[0]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]
}
